$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the existing header style (bold, bordered, centered) from H1
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$iValues = @(4,7,4,7,8,2,7,8,8,9,6,7,6,7,8,7,8,5,1,7,8,8)
$jValues = @(4,7,6,8,8,3,7,8,8,9,6,8,7,7,8,7,8,5,1,7,8,8)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
